$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (DATAEVENTO) held numeric Excel date/time serials; replace them
# with their formatted text equivalents "dd-MM-yyyy:HH:mm:ss".
$ws.Range("F2").Value = "16-07-2024:16:51:20"
$ws.Range("F3").Value = "16-07-2024:16:50:04"
$ws.Range("F4").Value = "16-07-2024:16:50:36"
$ws.Range("F5").Value = "16-07-2024:16:51:11"
$ws.Range("F6").Value = "16-07-2024:16:51:41"
$ws.Range("F7").Value = "16-07-2024:16:51:44"
